$wb = $excel.ActiveWorkbook

# --- Update the "Status" value for every locale row -------------------
# The localization run moved on from handoff into active translation, so
# every cell that used to read "Ready for handoff" now reads
# "In Translation" (Overview!E2/F2 for zh-cn/de-de, and the per-locale
# "Status" column C2 on the zh-cn and de-de sheets).

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the locale-status columns ----------------------------------
# The "zh-cn"/"de-de" columns on Overview, and the "Status" column on
# each locale sheet, get a bit narrower now that the report is archived.

$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E ("zh-cn")
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F ("de-de")

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5        # column C ("Status")
